$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("G8").Value = "dnasr281@gmail.com, System"
Write-Host ("VAL: " + $ws.Range("G8").Text)
